# Adding a third block of columns ("M_PL" = profits) to the table, mirroring
# the existing "M_%cit" (B:I) and "M_ETR" (J:Q) blocks, now placed in columns
# R:Y.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 18  # column R
$lastCol  = 25  # column Y

# --- Header row 1: merge the new block first, then stamp the same
#     bold / centered / boxed formatting used by the other header cells,
#     one cell at a time, so every cell in the merge ends up on the same
#     style index as the existing headers. ---
$ws.Range("R1:Y1").Merge() | Out-Null
for ($col = $firstCol; $col -le $lastCol; $col++) {
    $c = $ws.Cells.Item(1, $col)
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108   # xlCenter
    $c.VerticalAlignment = -4160     # xlTop
    $c.Borders.LineStyle = 1         # xlContinuous
}
$ws.Cells.Item(1, $firstCol).Value = "M_PL"

# --- Header row 2: repeat the 8 sub-column labels used in the other blocks,
#     with the same formatting as the existing row-2 headers. ---
$labels = @("GFA - Sales", "GFA - Sales + Emp", "IMF - Sales", "IMF - Sales + Emp", `
            "OECD (20%) - Sales", "OECD (20%) - Sales + Emp", "OECD - Sales", "OECD - Sales + Emp")
for ($col = $firstCol; $col -le $lastCol; $col++) {
    $c = $ws.Cells.Item(2, $col)
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4160
    $c.Borders.LineStyle = 1
}
for ($i = 0; $i -lt $labels.Length; $i++) {
    $ws.Cells.Item(2, $firstCol + $i).Value = $labels[$i]
}

# --- New data rows (R:Y) ---
$row4 = @(958366954868, 958617846484, 956175459402, 956426351018, 1008209699708, 1008209699708, 1008209699708, 1008209699708)
for ($i = 0; $i -lt $row4.Length; $i++) { $ws.Cells.Item(4, $firstCol + $i).Value = $row4[$i] }

$row5 = @(3140810, 3140810, 3140810, 3140810, 3140810, 3140810, 3140810, 3140810)
for ($i = 0; $i -lt $row5.Length; $i++) { $ws.Cells.Item(5, $firstCol + $i).Value = $row5[$i] }

$row6 = @(2557154421, 16062039693, 2124037026, 17056843395, 21277927825, 21277927825, 21277927825, 21277927825)
for ($i = 0; $i -lt $row6.Length; $i++) { $ws.Cells.Item(6, $firstCol + $i).Value = $row6[$i] }

$row7 = @(11464456004, 868911225, 1750040641, 868911225)
for ($i = 0; $i -lt $row7.Length; $i++) { $ws.Cells.Item(7, $firstCol + $i).Value = $row7[$i] }

$row8 = @(41567527900, 47033047618, 41286747272, 47534570308, 49527932043, 49527932043, 49527932043, 49527932043)
for ($i = 0; $i -lt $row8.Length; $i++) { $ws.Cells.Item(8, $firstCol + $i).Value = $row8[$i] }
